{"js": "const pairs = [\n  [\"2024-02-06 Tuesday\", \"2024-02-07 Wednesday\"],\n  [\"164\u00d77=\", \"276\u00d77=\"],\n  [\"718\u00d75=\", \"856\u00d72=\"],\n  [\"969\u00d77=\", \"859\u00d75=\"],\n  [\"927\u00d74=\", \"437\u00d78=\"],\n  [\"493\u00d74=\", \"589\u00d79=\"],\n  [\"579\u00d79=\", \"181\u00d74=\"],\n  [\"936\u00d72=\", \"335\u00d73=\"],\n  [\"518\u00d75=\", \"485\u00d75=\"],\n  [\"929\u00d79=\", \"516\u00d78=\"],\n  [\"303\u00d76=\", \"473\u00d75=\"],\n  [\"757\u00d74=\", \"798\u00d74=\"],\n  [\"384\u00d73=\", \"191\u00d77=\"],\n  [\"210\u00d76=\", \"739\u00d76=\"],\n  [\"971\u00d75=\", \"882\u00d79=\"],\n  [\"220\u00d78=\", \"233\u00d73=\"],\n  [\"808\u00d72=\", \"792\u00d79=\"],\n  [\"747\u00d77=\", \"297\u00d72=\"],\n  [\"588\u00d76=\", \"858\u00d75=\"],\n  [\"470\u00d73=\", \"564\u00d73=\"],\n  [\"651\u00d78=\", \"246\u00d73=\"],\n  [\"776\u00d72=\", \"908\u00d75=\"],\n  [\"305\u00d72=\", \"320\u00d79=\"],\n  [\"262\u00d78=\", \"285\u00d76=\"],\n  [\"983\u00d76=\", \"984\u00d74=\"],\n  [\"189\u00d73=\", \"882\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-06 Tuesday\", \"2024-02-07 Wednesday\"),\n    @(\"164\u00d77=\", \"276\u00d77=\"),\n    @(\"718\u00d75=\", \"856\u00d72=\"),\n    @(\"969\u00d77=\", \"859\u00d75=\"),\n    @(\"927\u00d74=\", \"437\u00d78=\"),\n    @(\"493\u00d74=\", \"589\u00d79=\"),\n    @(\"579\u00d79=\", \"181\u00d74=\"),\n    @(\"936\u00d72=\", \"335\u00d73=\"),\n    @(\"518\u00d75=\", \"485\u00d75=\"),\n    @(\"929\u00d79=\", \"516\u00d78=\"),\n    @(\"303\u00d76=\", \"473\u00d75=\"),\n    @(\"757\u00d74=\", \"798\u00d74=\"),\n    @(\"384\u00d73=\", \"191\u00d77=\"),\n    @(\"210\u00d76=\", \"739\u00d76=\"),\n    @(\"971\u00d75=\", \"882\u00d79=\"),\n    @(\"220\u00d78=\", \"233\u00d73=\"),\n    @(\"808\u00d72=\", \"792\u00d79=\"),\n    @(\"747\u00d77=\", \"297\u00d72=\"),\n    @(\"588\u00d76=\", \"858\u00d75=\"),\n    @(\"470\u00d73=\", \"564\u00d73=\"),\n    @(\"651\u00d78=\", \"246\u00d73=\"),\n    @(\"776\u00d72=\", \"908\u00d75=\"),\n    @(\"305\u00d72=\", \"320\u00d79=\"),\n    @(\"262\u00d78=\", \"285\u00d76=\"),\n    @(\"983\u00d76=\", \"984\u00d74=\"),\n    @(\"189\u00d73=\", \"882\u00d72=\")\n)\n\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($before, $false, $false, $false, $false, $false, $true, 1, $false, $after, 2) | Out-Null\n}\n"}
